{"js": "// Replace the date and the 25 division-problem texts per the commit diff.\nconst replacements = [\n  [\"2025-01-30 Thursday\", \"2025-01-31 Friday\"],\n  [\"739\u00f76=\", \"455\u00f77=\"],\n  [\"504\u00f77=\", \"380\u00f79=\"],\n  [\"629\u00f79=\", \"634\u00f78=\"],\n  [\"309\u00f76=\", \"246\u00f79=\"],\n  [\"934\u00f74=\", \"290\u00f79=\"],\n  [\"486\u00f77=\", \"264\u00f74=\"],\n  [\"468\u00f75=\", \"643\u00f77=\"],\n  [\"988\u00f77=\", \"840\u00f74=\"],\n  [\"967\u00f77=\", \"195\u00f73=\"],\n  [\"764\u00f78=\", \"313\u00f73=\"],\n  [\"747\u00f75=\", \"601\u00f75=\"],\n  [\"908\u00f75=\", \"370\u00f76=\"],\n  [\"277\u00f77=\", \"858\u00f75=\"],\n  [\"391\u00f75=\", \"407\u00f78=\"],\n  [\"361\u00f76=\", \"316\u00f79=\"],\n  [\"638\u00f76=\", \"473\u00f76=\"],\n  [\"795\u00f75=\", \"544\u00f72=\"],\n  [\"433\u00f76=\", \"568\u00f76=\"],\n  [\"420\u00f72=\", \"930\u00f72=\"],\n  [\"664\u00f78=\", \"467\u00f76=\"],\n  [\"777\u00f75=\", \"440\u00f75=\"],\n  [\"301\u00f72=\", \"607\u00f72=\"],\n  [\"282\u00f77=\", \"636\u00f74=\"],\n  [\"913\u00f76=\", \"417\u00f78=\"],\n  [\"303\u00f75=\", \"631\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and the 25 division-problem texts per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-01-30 Thursday\", \"2025-01-31 Friday\"),\n  @(\"739\u00f76=\", \"455\u00f77=\"),\n  @(\"504\u00f77=\", \"380\u00f79=\"),\n  @(\"629\u00f79=\", \"634\u00f78=\"),\n  @(\"309\u00f76=\", \"246\u00f79=\"),\n  @(\"934\u00f74=\", \"290\u00f79=\"),\n  @(\"486\u00f77=\", \"264\u00f74=\"),\n  @(\"468\u00f75=\", \"643\u00f77=\"),\n  @(\"988\u00f77=\", \"840\u00f74=\"),\n  @(\"967\u00f77=\", \"195\u00f73=\"),\n  @(\"764\u00f78=\", \"313\u00f73=\"),\n  @(\"747\u00f75=\", \"601\u00f75=\"),\n  @(\"908\u00f75=\", \"370\u00f76=\"),\n  @(\"277\u00f77=\", \"858\u00f75=\"),\n  @(\"391\u00f75=\", \"407\u00f78=\"),\n  @(\"361\u00f76=\", \"316\u00f79=\"),\n  @(\"638\u00f76=\", \"473\u00f76=\"),\n  @(\"795\u00f75=\", \"544\u00f72=\"),\n  @(\"433\u00f76=\", \"568\u00f76=\"),\n  @(\"420\u00f72=\", \"930\u00f72=\"),\n  @(\"664\u00f78=\", \"467\u00f76=\"),\n  @(\"777\u00f75=\", \"440\u00f75=\"),\n  @(\"301\u00f72=\", \"607\u00f72=\"),\n  @(\"282\u00f77=\", \"636\u00f74=\"),\n  @(\"913\u00f76=\", \"417\u00f78=\"),\n  @(\"303\u00f75=\", \"631\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
